$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Designated Team names were renamed:
#   Team1 -> Jason
#   Team2 -> TST
#   Team3 -> max
# and the team assignment for rows 11 and 13 was swapped at the same time
# (row 11 used to be Team3, row 13 used to be Team1).
$ws.Range("U3").Value = "Jason"
$ws.Range("U4").Value = "TST"
$ws.Range("U5").Value = "max"
$ws.Range("U11").Value = "Jason"
$ws.Range("U12").Value = "TST"
$ws.Range("U13").Value = "max"

# Merge the split "R11" + "R12:R13" text-length validation into one
# contiguous "R11:R13" range (and move it after the B11:B13 / C11:C13
# rules, matching the re-saved rule order).
$ws.Range("R11").Validation.Delete() | Out-Null
$ws.Range("R12:R13").Validation.Delete() | Out-Null
$ws.Range("R11:R13").Validation.Add(6, 1, 1, "6", "6") | Out-Null

# Active cell/selection moved from T23 to T22
$ws.Range("T22").Select() | Out-Null
